$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Date: bump the generation timestamp
$meta.Range("B8").Value = "2026-01-01T13:37:23+00:00"

# Description: drop the trailing "Supports goal-directed..." sentence
$meta.Range("B11").Value = "Extension to link nursing interventions to the patient goals they are intended to achieve."

# --- Elements sheet updates ---
$elems = $wb.Worksheets.Item("Elements")

# The root Extension row's Definition repeats the same description text, keep it in sync
$elems.Range("M2").Value = "Extension to link nursing interventions to the patient goals they are intended to achieve."

# Extension.value[x] Type(s): point the reference at the renamed goal profile,
# keep the trailing newline that was already part of the cell text
$elems.Range("K6").Value = "Reference(https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/onc-nursing-goal)`n"
